# Auto-generated update of cryptos worksheet to reflect latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.161.70'
$ws.Range("E2").Value = '  +1.42%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.422.38'
$ws.Range("E3").Value = '  +1.67%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '558.75'
$ws.Range("E5").Value = '  +1.59%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.61'
$ws.Range("E6").Value = '  +3.15%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.532'
$ws.Range("E8").Value = '  +1.36%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.420.07'
$ws.Range("E9").Value = '  +1.46%  '

# Row 10
$ws.Range("E10").Value = '  +0.33%  '

# Row 11
$ws.Range("E11").Value = '  -1.45%  '

# Row 12
$ws.Range("E12").Value = '  +1.08%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.353'
$ws.Range("E13").Value = '  +0.79%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.22'
$ws.Range("E14").Value = '  +4.26%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000176'
$ws.Range("E15").Value = '  +5.44%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.848.53'
$ws.Range("E16").Value = '  +1.88%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.949.69'
$ws.Range("E17").Value = '  +1.19%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.417.47'
$ws.Range("E18").Value = '  +1.48%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.19'
$ws.Range("E19").Value = '  +2.36%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.19'
$ws.Range("E20").Value = '  +0.66%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '324.15'
$ws.Range("E21").Value = '  +0.50%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.77'
$ws.Range("E22").Value = '  +0.36%  '

# Row 23
$ws.Range("E23").Value = '  -0.03%  '

# Row 24
$ws.Range("E24").Value = '  +1.38%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.72'
$ws.Range("E25").Value = '  +0.43%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.98'
$ws.Range("E26").Value = '  +5.18%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '594.73'
$ws.Range("E27").Value = '  +15.38%  '

# Row 28
$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.536.43'
$ws.Range("E28").Value = '  +1.88%  '

# Row 29
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.30%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0940'
$ws.Range("E30").Value = '  +4.84%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.30'
$ws.Range("E31").Value = '  +1.10%  '

# Row 32
$ws.Range("E32").Value = '  +4.33%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.150'
$ws.Range("E33").Value = '  -0.34%  '

# Row 34
$ws.Range("E34").Value = '  +1.87%  '

# Row 35
$ws.Range("E35").Value = '  +1.65%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.71'
$ws.Range("E36").Value = '  +4.57%  '

# Row 37
$ws.Range("E37").Value = '  +0.09%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.78'
$ws.Range("E38").Value = '  +1.50%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.385'
$ws.Range("E39").Value = '  +1.38%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '152.22'
$ws.Range("E40").Value = '  +3.81%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.72'
$ws.Range("E41").Value = '  +0.29%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.82'
$ws.Range("E42").Value = '  -4.95%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.37'
$ws.Range("E44").Value = '  +10.35%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '150.96'
$ws.Range("E45").Value = '  -0.56%  '

# Row 46
$ws.Range("E46").Value = '  +1.02%  '

# Row 47
$ws.Range("E47").Value = '  +2.75%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.27'
$ws.Range("E48").Value = '  +4.41%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.592'
$ws.Range("E49").Value = '  +2.32%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0920'
$ws.Range("E50").Value = '  +0.91%  '

# Row 51
$ws.Range("E51").Value = '  +1.93%  '
